$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at N, shifting old N/O/P to O/P/Q
$ws.Columns("N").Insert()

# Set the width of the newly inserted column N (non-bestFit, custom width 10)
$ws.Columns("N").ColumnWidth = 9.1666667

# Update the active cell selection to match the new editing position
[void]$ws.Range("R7").Select()
